# Daily Satellite Data Update
# Updates the Prelety sheet with refreshed flyover timing data for the
# existing date rows (2-8) and appends two new rows (9-10) for 2026-02-27.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-BGR($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $b*65536 + $g*256 + $r
}

# ---------------------------------------------------------------------
# Row data for columns A..R (rows 2 through 10)
# ---------------------------------------------------------------------
$data = @(
  @("20260221--01",14,"04:01","00:00","05:07:12","05:10:04","05:12:05","05:14:05","05:16:58","-","-",-9.1,"A+B","2",96,92,96,96),
  @("20260222--01",13,"03:08","03:08","04:28:54","04:32:03","04:33:37","04:35:11","04:38:20","4°","04:30:11",-15,"A","4",100,81,100,100),
  @("20260223--01",11,"01:34","01:34","03:50:36","03:54:20","03:55:06","03:55:54","03:59:38","9°","03:53:37",-20.7,"A","4",93,80,53,45),
  @("20260224--01",24,"05:51","05:51","04:48:25","04:50:49","04:53:44","04:56:40","04:59:04","8°","04:50:26",-11.2,"A","2",95,93,9,0),
  @("20260225--01",22,"05:36","04:20","04:09:51","04:12:18","04:15:06","04:17:54","04:20:22","17°","04:13:34",-17,"A+B","2",58,29,40,18),
  @("20260226--01",20,"05:18","02:29","03:31:15","03:33:46","03:36:25","03:39:04","03:41:36","20°","03:36:35",-22.7,"A+B","2",37,37,0,7),
  @("20260226--02",32,"06:27","06:27","05:08:22","05:10:38","05:13:51","05:17:05","05:19:21","7°","05:10:07",-7.4,"A","1",42,42,0,7),
  @("20260227--01",18,"04:52","00:38","02:52:38","02:55:16","02:57:42","03:00:08","03:02:46","13°","02:59:30",-28.1,"A+B","2",94,28,84,92),
  @("20260227--02",32,"06:24","05:16","04:29:37","04:31:54","04:35:06","04:38:18","04:40:35","18°","04:33:02",-13.2,"A+B","2",83,22,83,83)
)

# Fill colors (hex, no alpha) for columns O,P,Q,R for rows 2..10
$fillColors = @(
  @("F0707F","EEF3F9","F6F9FC","F6F9FC"),
  @("F0707F","DDE8F4","FFFFFF","FFFFFF"),
  @("F0707F","DDE8F4","B2CBE6","A2C0E1"),
  @("F0707F","F6F9FC","6697CD","558CC8"),
  @("F0707F","88AED8","99BADE","77A3D3"),
  @("F0B070","90B4DB","558CC8","5E92CB"),
  @("F0707F","99BADE","558CC8","5E92CB"),
  @("F0707F","88AED8","E6EEF7","EEF3F9"),
  @("F0707F","77A3D3","E6EEF7","E6EEF7")
)

$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")
$fontGray = Get-BGR("333333")

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]

    for ($c = 0; $c -lt $columns.Count; $c++) {
        $colLetter = $columns[$c]
        $cell = $ws.Range($colLetter + $rowNum)
        $cell.Value = $rowVals[$c]
    }

    # Column N holds small integer "priority" codes that must remain TEXT
    # (as in the source file) rather than become numeric. Re-apply the
    # formatting of a known-text sibling cell (column A, same style) so the
    # numeric-looking value is forced back to a shared string without
    # disturbing the cell's visual style.
    $nCell = $ws.Range("N" + $rowNum)
    $nCell.NumberFormat = "@"
    $nCell.Value = [string]$rowVals[13]
    $ws.Range("A2").Copy() | Out-Null
    $nCell.PasteSpecial(-4122) | Out-Null  # xlPasteFormats

    # Apply the heat-map style colors to O,P,Q,R
    $colors = $fillColors[$i]
    $oprColumns = @("O","P","Q","R")
    for ($k = 0; $k -lt 4; $k++) {
        $target = $ws.Range($oprColumns[$k] + $rowNum)
        $target.Interior.Color = Get-BGR($colors[$k])
        $target.Font.Color = $fontGray
        $target.HorizontalAlignment = -4108  # xlCenter
        $target.VerticalAlignment = -4108    # xlCenter
    }
}

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Extend conditional formatting ranges from row 8 to row 10 for all
# data columns (A:R) so the new rows participate in the existing rules.
# ---------------------------------------------------------------------
foreach ($colLetter in $columns) {
    $oldRange = $ws.Range($colLetter + "2:" + $colLetter + "8")
    $newRange = $ws.Range($colLetter + "2:" + $colLetter + "10")
    $fcs = $oldRange.FormatConditions
    $cnt = $fcs.Count()
    for ($j = 1; $j -le $cnt; $j++) {
        $fc = $fcs.Item($j)
        $fc.ModifyAppliesToRange($newRange) | Out-Null
    }
}
